# The deck ships two DrawingML theme parts:
#   ppt/theme/theme1.xml  -> was "Office Theme"  (only used by the Notes Master)
#   ppt/theme/theme2.xml  -> was "Integral"       (used by the one real Slide Master,
#                                                   i.e. by every slide/layout)
#
# The authored edit swaps the two themes' contents, so the Slide Master (and
# therefore every slide) ends up painted with the plain "Office Theme" palette
# instead of "Integral", while the Notes Master swaps the other way.
#
# The PowerPoint object model has no "replace this theme part wholesale"
# call, so we reproduce the edit the way a real user would in the UI: by
# recolouring the active theme's 12 scheme colors (Background/Text 1-2,
# Accent 1-6, Hyperlink, Followed Hyperlink) through
# Theme.ThemeColorScheme.Colors(i).RGB - the exact values "Office Theme"
# used for its <a:clrScheme>.

$p  = $ppt.ActivePresentation
$sm = $p.SlideMaster
$tcs = $sm.Theme.ThemeColorScheme

function Set-SchemeColor {
    param(
        [int]$Index,
        [string]$Hex
    )
    $r = [Convert]::ToInt32($Hex.Substring(0,2), 16)
    $g = [Convert]::ToInt32($Hex.Substring(2,2), 16)
    $b = [Convert]::ToInt32($Hex.Substring(4,2), 16)
    # VBA's RGB()/ColorFormat.RGB long uses 0x00BBGGRR ordering.
    $tcs.Colors($Index).RGB = $r + ($g * 256) + ($b * 65536)
}

# MsoThemeColorSchemeIndex order: 1 dk1, 2 lt1, 3 dk2, 4 lt2,
# 5-10 accent1-6, 11 hlink, 12 folHlink — the "Office Theme" values.
Set-SchemeColor 1  "000000"
Set-SchemeColor 2  "FFFFFF"
Set-SchemeColor 3  "44546A"
Set-SchemeColor 4  "E7E6E6"
Set-SchemeColor 5  "5B9BD5"
Set-SchemeColor 6  "ED7D31"
Set-SchemeColor 7  "A5A5A5"
Set-SchemeColor 8  "FFC000"
Set-SchemeColor 9  "4472C4"
Set-SchemeColor 10 "70AD47"
Set-SchemeColor 11 "0563C1"
Set-SchemeColor 12 "954F72"
